# ----------------------------------------------------------------------------
# "Mise a jour site" - price update for prix/bicones.xlsx
#
# - Column D (prix_tvac) gets a ~21% price increase on every data row (2-72),
#   matching the new price list.
# - Column E (Article) is highlighted to flag the rows whose price changed:
#   yellow for the bulk of the updated articles (rows 2-70) and green for the
#   two newly added/adjusted "Culasse" articles at the bottom (rows 71-72).
# - The last action on the sheet leaves the selection on D73 (just below the
#   last data row), matching the editor state when the workbook was saved.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update prices (column D, prix_tvac) ---
$ws.Range("D2").Value = 4.92
$ws.Range("D3").Value = 5.3100000000000005
$ws.Range("D4").Value = 4.86
$ws.Range("D5").Value = 4.89
$ws.Range("D6").Value = 5.03
$ws.Range("D7").Value = 7.82
$ws.Range("D8").Value = 8.17
$ws.Range("D9").Value = 13.1
$ws.Range("D10").Value = 6.49
$ws.Range("D11").Value = 4.86
$ws.Range("D12").Value = 6.28
$ws.Range("D13").Value = 3.11
$ws.Range("D14").Value = 4.62
$ws.Range("D15").Value = 4.62
$ws.Range("D16").Value = 13.25
$ws.Range("D17").Value = 5.45
$ws.Range("D18").Value = 6.3
$ws.Range("D19").Value = 8.25
$ws.Range("D20").Value = 11.700000000000001
$ws.Range("D21").Value = 7.18
$ws.Range("D22").Value = 5.9
$ws.Range("D23").Value = 6.99
$ws.Range("D24").Value = 7.68
$ws.Range("D25").Value = 7.38
$ws.Range("D26").Value = 12.5
$ws.Range("D27").Value = 9.56
$ws.Range("D28").Value = 24.27
$ws.Range("D29").Value = 9.05
$ws.Range("D30").Value = 6.63
$ws.Range("D31").Value = 4.16
$ws.Range("D32").Value = 5.01
$ws.Range("D33").Value = 6.3500000000000005
$ws.Range("D34").Value = 7.3500000000000005
$ws.Range("D35").Value = 6.12
$ws.Range("D36").Value = 6.36
$ws.Range("D37").Value = 4.62
$ws.Range("D38").Value = 5.4
$ws.Range("D39").Value = 6.12
$ws.Range("D40").Value = 6.12
$ws.Range("D41").Value = 10.07
$ws.Range("D42").Value = 4.6000000000000005
$ws.Range("D43").Value = 3.94
$ws.Range("D44").Value = 4.28
$ws.Range("D45").Value = 4.6000000000000005
$ws.Range("D46").Value = 5.3
$ws.Range("D47").Value = 4.1
$ws.Range("D48").Value = 4.95
$ws.Range("D49").Value = 3.46
$ws.Range("D50").Value = 4.04
$ws.Range("D51").Value = 3.74
$ws.Range("D52").Value = 5.63
$ws.Range("D53").Value = 6.86
$ws.Range("D54").Value = 5.45
$ws.Range("D55").Value = 11.3
$ws.Range("D56").Value = 12.05
$ws.Range("D57").Value = 3.11
$ws.Range("D58").Value = 3.04
$ws.Range("D59").Value = 3.34
$ws.Range("D60").Value = 4.21
$ws.Range("D61").Value = 4.05
$ws.Range("D62").Value = 5.01
$ws.Range("D63").Value = 3.9
$ws.Range("D64").Value = 4.04
$ws.Range("D65").Value = 4.04
$ws.Range("D66").Value = 3.86
$ws.Range("D67").Value = 5.8
$ws.Range("D68").Value = 6.86
$ws.Range("D69").Value = 6.2700000000000005
$ws.Range("D70").Value = 12.790000000000001
$ws.Range("D71").Value = 7.6
$ws.Range("D72").Value = 6.03

# --- Highlight the updated article codes (column E) ---
# Yellow fill (RGB 255,255,0) for rows 2-70
$ws.Range("E2:E70").Interior.Color = 65535
# Green fill (RGB 146,208,80) for rows 71-72
$ws.Range("E71:E72").Interior.Color = 5296274

# --- Restore the cursor/selection position used when the file was saved ---
$ws.Range("D73").Select()
